$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted above the current row 169,
# pushing the existing rows 169-238 down to 170-239. Excel's native
# row-insert shifts all the cell data (and styles) down for us.
$ws.Rows.Item(169).Insert()

# Populate the newly-inserted row 169 with the new observation. All
# columns match the record that used to sit at row 169 (now at row 170)
# except the date (D) and volume (J), which carry the new values.
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 44704
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 100112009
$ws.Cells.Item(169, 7).Value = "Acelga"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 200
$ws.Cells.Item(169, 11).Value = 550
$ws.Cells.Item(169, 12).Value = 600
$ws.Cells.Item(169, 13).Value = 575
$ws.Cells.Item(169, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(169, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(169, 16).Value = 575
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
